# Insert a new record row just before the current row 33 (case 5940 /
# SANCHEZ DE LORIA 1406). This shifts every row from the old 33..53 down
# by one (new 34..54), matching the diff exactly (old row 53, case 6137,
# ends up at new row 54; dimension grows from N53 to N54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 33 - pushes existing rows 33..53 down to 34..54
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record.
# Leading "'" forces text storage for numeric-looking values (Caso, Comuna,
# OT, Attachments) so they match the existing column typing (plain text),
# instead of being auto-coerced to numbers/dates by Excel.
$ws.Range("A33").Value = "'5883"
$ws.Range("B33").Value = "'5/26/2025"
$ws.Range("C33").Value = "CONGRESO AV. 2699"
$ws.Range("D33").Value = "'13"
$ws.Range("E33").Value = "'806944763"
$ws.Range("F33").Value = "PEBCOM"
$ws.Range("G33").Value = "Pendiente"
$ws.Range("H33").Value = "Picada"
$ws.Range("I33").Value = "'1"
$ws.Range("J33").Value = "Cambio"
$ws.Range("K33").Value = "Sin equipos"
$ws.Range("L33").Value = "Pasante"
$ws.Range("M33").Value = -58.46522
$ws.Range("N33").Value = -34.556786

# Clear the auto-applied "Text" number-format style picked up from the
# leading-apostrophe text coercion above, so the new cells carry no direct
# formatting, matching the plain (unstyled) data rows around them.
$ws.Range("A33:N33").Style = "Normal"
